# Reshape the results table:
# Old layout: A,B,C(pearson_r),D(pearson_p),E(spearman_r),F(spearman_p),G(n),H(mean_diff)
# New layout: A,B,C(spearman_r),D(spearman_p),E(n),F(mean_diff [new values])

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row
$ws.Range("C1").Value = "spearman_r"
$ws.Range("D1").Value = "spearman_p"
$ws.Range("E1").Value = "n"
$ws.Range("F1").Value = "mean_diff"

# Clear old G and H columns (no longer part of the table)
$ws.Range("G1:H9").Clear()

# New data values, per row: C=spearman_r, D=spearman_p, E=n, F=mean_diff(new)
$ws.Cells.Item(2, 3).Value = [double]"0.02578971968863382"
$ws.Cells.Item(2, 4).Value = [double]"0.6437244148461119"
$ws.Cells.Item(2, 5).Value = [double]"324"
$ws.Cells.Item(2, 6).Value = [double]"92.25689300411524"

$ws.Cells.Item(3, 3).Value = [double]"-0.5366818968278821"
$ws.Cells.Item(3, 4).Value = [double]"1.44973140734144e-25"
$ws.Cells.Item(3, 5).Value = [double]"324"
$ws.Cells.Item(3, 6).Value = [double]"112.2081275720165"

$ws.Cells.Item(4, 3).Value = [double]"0.07791892983253103"
$ws.Cells.Item(4, 4).Value = [double]"0.1617385594216551"
$ws.Cells.Item(4, 5).Value = [double]"324"
$ws.Cells.Item(4, 6).Value = [double]"1.134897119341559"

$ws.Cells.Item(5, 3).Value = [double]"-0.5386697465698866"
$ws.Cells.Item(5, 4).Value = [double]"8.90155837032734e-26"
$ws.Cells.Item(5, 5).Value = [double]"324"
$ws.Cells.Item(5, 6).Value = [double]"21.08613168724279"

$ws.Cells.Item(6, 3).Value = [double]"0.1600158738541625"
$ws.Cells.Item(6, 4).Value = [double]"0.003879661697487937"
$ws.Cells.Item(6, 5).Value = [double]"324"
$ws.Cells.Item(6, 6).Value = [double]"107.7940118312757"

$ws.Cells.Item(7, 3).Value = [double]"0.1369346169458027"
$ws.Cells.Item(7, 4).Value = [double]"0.01362777088369083"
$ws.Cells.Item(7, 5).Value = [double]"324"
$ws.Cells.Item(7, 6).Value = [double]"127.7452463991769"

$ws.Cells.Item(8, 3).Value = [double]"0.08165690046153645"
$ws.Cells.Item(8, 4).Value = [double]"0.1424866265162245"
$ws.Cells.Item(8, 5).Value = [double]"324"
$ws.Cells.Item(8, 6).Value = [double]"12145.06078549383"

$ws.Cells.Item(9, 3).Value = [double]"-0.5053471280814221"
$ws.Cells.Item(9, 4).Value = [double]"2.099657985557351e-22"
$ws.Cells.Item(9, 5).Value = [double]"324"
$ws.Cells.Item(9, 6).Value = [double]"12165.01202006173"
